# Apply the target edits to the workbook:
#  - Summary sheet: Fee Drag % (B6) goes from 3.85 to 0
#  - Strategies sheet: row 2 changes from a "Short Call" trade held 4 days
#    (2 in fees, 12.5 theta/day) to a "Short Put" trade held 2 hours
#    (0 in fees, 600 theta/day), with the segments JSON text updated to match.

$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 0

$wsStrategies = $wb.Worksheets.Item("Strategies")
$wsStrategies.Range("C2").Value = "Short Put"
$wsStrategies.Range("E2").Value = 0
$wsStrategies.Range("F2").Value = 0.08333333333333333
$wsStrategies.Range("G2").Value = 600
$wsStrategies.Range("I2").Value = "[{'strategy_name': 'Short Put', 'pnl': 50.0, 'entry_ts': '2025-01-01T10:00:00', 'exit_ts': '2025-01-01T12:00:00'}]"
